$d = $word.ActiveDocument

# --- 1) "In the case of the MICE switches at 511, Cologix accepts ..."
#        -> "In the case of the MICE switches, Cologix and Ridgeview accept ..."
$d.Content.Find.Execute("MICE switches at 511, ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "MICE switches, ", 2) | Out-Null

$d.Content.Find.Execute("Cologix accepts", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Cologix and Ridgeview accept", 2) | Out-Null

# --- 2) "... have a Cologix account representative, contact Scott McCrady <...>."
#        -> "... have an account representative… For Cologix, contact Scott McCrady <...>."
$d.Content.Find.Execute("have a Cologix account representative, contact", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ("have an account representative" + [char]0x2026 + " For Cologix, contact"), 2) | Out-Null

# --- 3) Append new sentence about Ridgeview right after the existing Cologix sentence.
#        The existing sentence ends in the unique literal ">." immediately followed by a
#        new paragraph, so anchor on that (it does not overlap the hyperlink run).
$tail = $d.Content
$tail.Find.Execute(">.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail.Collapse(0)
$tail.InsertAfter(" For Ridgeview, contact Michael Dumas <")

$lead = $d.Content
$lead.Find.Execute(" For Ridgeview, contact Michael Dumas <", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$lead.Font.Name = "Arial"
$lead.Font.NameBi = "Arial"
$lead.Collapse(0)
$lead.InsertAfter("MDumas@usinternet.com")

$linkRange = $d.Content
$linkRange.Find.Execute("MDumas@usinternet.com", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($linkRange, "mailto:MDumas@usinternet.com", "", "", "MDumas@usinternet.com") | Out-Null

$linkFont = $d.Content
$linkFont.Find.Execute("MDumas@usinternet.com", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$linkFont.Font.Name = "Arial"
$linkFont.Collapse(0)
$linkFont.InsertAfter(">.")

$closingFont = $d.Content
$closingFont.Find.Execute("MDumas@usinternet.com>.", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0) | Out-Null
